$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 for the newly published dissertation-based
# journal article (it previously lived as an "under review" entry further
# down the sheet; that duplicate entry is removed below).
$ws.Rows(2).Insert()

$ws.Range("A2").Value = "prpa "
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Wright, J., Sohlberg, M.M., McIntosh, K., Seeley, J., Hadley, W., Blitz, D. & Lowham, E."
$ws.Range("D2").Value = 2022
$ws.Range("E2").Value = "What is the effect of personalized cognitive strategy instruction on facilitating return-to-learn for individuals experiencing prolonged concussion symptoms"
$ws.Range("F2").Value = "Neuropsychological Rehabilitation "
$ws.Range("J2").Value = "https://doi.org/10.1080/09602011.2022.2074467"
$ws.Range("K2").Value = "https://www.tandfonline.com/doi/full/10.1080/09602011.2022.2074467"
$ws.Range("M2").Value = "Department of Communication Disorders & Sciences, University of Oregon"

# Clear the columns that this row doesn't use (Insert() copies formatting
# from the row above across the full width; drop the untouched cells so
# the row only carries the columns it actually needs).
$ws.Range("G2:I2").Clear()
$ws.Range("L2").Clear()
$ws.Range("N2:W2").Clear()

$ws.Rows(2).RowHeight = 119

# Row 5 (the "under review" duplicate of the dissertation study, now
# superseded by the published row above) is removed entirely.
$ws.Rows(5).Delete()

# Row 6 (the unpublished Hawaii retrospective manuscript) is also removed.
$ws.Rows(6).Delete()

$ws.Range("A1").Select()
